$d = $word.ActiveDocument

# The sentencing/suspension date moves from June 24, 2022 -> June 26, 2022.
# This exact string occurs 3 times in the document (the narrative paragraph,
# the bold "License Suspension" date field, and the suspension clause);
# Find/Execute with Replace:=wdReplaceAll updates every occurrence.
$d.Content.Find.Execute("June 24, 2022", $true, $false, $false, $false, $false, `
    $true, 1, $false, "June 26, 2022", 2)

# The community-service proof-of-completion deadline moves from
# August 23, 2022 -> August 25, 2022.
$d.Content.Find.Execute("August 23, 2022", $true, $false, $false, $false, $false, `
    $true, 1, $false, "August 25, 2022", 2)
